$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (measure/dimension type row)
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "iaest-dimension:superficie-util"

# Row 3 (medida/dim row)
$ws.Range("C3").Value = "dim"
$ws.Range("G3").Value = "dim"

# Row 4 (type/URI row)
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "skos:Concept"

# Row 5 (mapping file row) - move value from E5 to G5 with new name,
# carrying over the same cell formatting style.
$ws.Range("E5").Copy()
$ws.Range("G5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E5").Clear()
$ws.Range("G5").Value = "mapping-superficie-util.xlsx"
